$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 13, shifting existing rows 13:112 down to 14:113
$ws.Rows("13:13").Insert()

# Populate the newly inserted row 13 with the new data sample
$ws.Range("A13").Value = 11
$ws.Range("B13").Value = "Vega Monumental Concepción"
$ws.Range("C13").Value = "Bíobío"
$ws.Range("D13").Value = (Get-Date -Year 2023 -Month 4 -Day 6 -Hour 0 -Minute 0 -Second 0)
$ws.Range("E13").Value = 8
$ws.Range("F13").Value = 100112012
$ws.Range("G13").Value = "Espinaca"
$ws.Range("H13").Value = "Sin especificar"
$ws.Range("I13").Value = "Primera"
$ws.Range("J13").Value = 40
$ws.Range("K13").Value = 9000
$ws.Range("L13").Value = 10000
$ws.Range("M13").Value = 9500
$ws.Range("N13").Value = "`$/cuna 10 kilos"
$ws.Range("O13").Value = "Región Metropolitana"
$ws.Range("P13").Value = 950
$ws.Range("Q13").Value = 10
$ws.Range("R13").Value = "Hortaliza"
